$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A14").Value = 111251432
$ws.Range("B14").Value = 81248
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 1312
$ws.Range("F14").Value = "Gammelgransskål"
$ws.Range("G14").Value = "Pseudographis pinicola"
$ws.Range("H14").Value = "(Nyl.) Rehm"
$ws.Range("Q14").Value = 460622.5513675315
$ws.Range("R14").Value = 7165027.330594921

$ws.Range("A15").Value = 111251407
$ws.Range("B15").Value = 73696
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 6440
$ws.Range("F15").Value = "Vitgrynig nållav"
$ws.Range("G15").Value = "Chaenotheca subroscida"
$ws.Range("H15").Value = "(Eitner) Zahlbr."
$ws.Range("Q15").Value = 460240.5118381025
$ws.Range("R15").Value = 7164805.620072429

$ws.Range("A16").Value = 111251402
$ws.Range("B16").Value = 73696
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6440
$ws.Range("F16").Value = "Vitgrynig nållav"
$ws.Range("G16").Value = "Chaenotheca subroscida"
$ws.Range("H16").Value = "(Eitner) Zahlbr."
$ws.Range("Q16").Value = 460212.3128264685
$ws.Range("R16").Value = 7164818.870384302

$ws.Range("A17").Value = 111251434
$ws.Range("B17").Value = 78612
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 6464
$ws.Range("F17").Value = "Luddlav"
$ws.Range("G17").Value = "Nephroma resupinatum"
$ws.Range("H17").Value = "(L.) Ach."
$ws.Range("Q17").Value = 460452.9763639791
$ws.Range("R17").Value = 7164846.208533676

$ws.Range("A18").Value = 111251423
$ws.Range("B18").Value = 77677
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 1249
$ws.Range("F18").Value = "Norsk näverlav"
$ws.Range("G18").Value = "Platismatia norvegica"
$ws.Range("H18").Value = "(Lynge) W.L.Culb. & C.F.Culb."
$ws.Range("Q18").Value = 460188.7895233887
$ws.Range("R18").Value = 7164860.82616597

$ws.Range("A19").Value = 111251428
$ws.Range("B19").Value = 89423
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 5432
$ws.Range("F19").Value = "Granticka"
$ws.Range("G19").Value = "Porodaedalea chrysoloma"
$ws.Range("H19").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q19").Value = 460445.0942901828
$ws.Range("R19").Value = 7164835.148113105

$ws.Range("A20").Value = 111251430
$ws.Range("B20").Value = 77515
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("Q20").Value = 460188.8289468794
$ws.Range("R20").Value = 7164863.831099218

$ws.Range("A22").Value = 111251420
$ws.Range("B22").Value = 77677
$ws.Range("D22").Value = "VU"
$ws.Range("E22").Value = 1249
$ws.Range("F22").Value = "Norsk näverlav"
$ws.Range("G22").Value = "Platismatia norvegica"
$ws.Range("H22").Value = "(Lynge) W.L.Culb. & C.F.Culb."
$ws.Range("Q22").Value = 460243.4530616797
$ws.Range("R22").Value = 7164800.429238674
